# ---------------------------------------------------------------------------
# "expansao das analises automaticas"
#
# 1) particip (E) and taxa_sucesso (F) are rescaled from fractions (0..1) to
#    percentage points (0..100) -- i.e. every existing value in those two
#    columns is multiplied by 100.
# 2) Three new metric columns are appended after max_sucesso (K):
#       L  apoio_medio
#       M  contribuicoes
#       N  media_contribuicoes
#    with their own header cells on row 1 and a value on every data row
#    (rows 2-74), matching the header's look (font/border/alignment).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells (L1:N1) ------------------------------------------
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Give the 3 new header cells the same look as the existing header row
# (bold font + border + centered alignment) by copying K1's format onto them.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Per-row data --------------------------------------------------------
# New E (particip) / F (taxa_sucesso) values -- x100 of the original -- plus
# the new L (apoio_medio) / M (contribuicoes) / N (media_contribuicoes)
# values, one CSV line per worksheet row, rows 2..74 in order.
$data = @"
0.0749063670411985,0,0,0,0
0.5243445692883895,71.42857142857143,72.63883144032484,769,153.8
1.048689138576779,28.57142857142857,69.16365343052672,495,123.75
1.423220973782771,36.84210526315789,95.50962732859578,1073,153.2857142857143
1.573033707865169,42.85714285714285,71.38657292738846,1618,179.7777777777778
2.9962546816479403,57.49999999999999,92.98191948624745,5497,239
0.8239700374531835,36.36363636363636,112.2443057692103,649,162.25
1.048689138576779,71.42857142857143,55.85601551500451,1931,193.1
0.22471910112359547,33.33333333333333,52.78244028225377,37,37
8.614232209737828,58.26086956521739,84.51296022186848,21560,321.7910447761194
0.149812734082397,50.0,71.24263946639917,544,544
0.149812734082397,100,67.93158644297183,239,119.5
0.8988764044943821,33.33333333333333,72.42130976357548,1698,424.5
1.4981273408239701,50.0,74.91086123718691,3894,389.4
3.071161048689139,58.536585365853654,81.12766312510587,4348,181.1666666666667
1.4981273408239701,65.0,72.24040738498717,3100,238.4615384615385
6.067415730337079,59.25925925925925,81.23123121034725,22855,476.1458333333333
11.38576779026217,61.8421052631579,86.86464629097692,22939,244.031914893617
1.4981273408239701,55.00000000000001,92.89347465855401,1717,156.0909090909091
0.0749063670411985,0,0,0,0
7.640449438202247,65.68627450980392,82.90726758980419,18775,280.2238805970149
1.947565543071161,53.84615384615385,106.1211981629064,3590,256.4285714285714
0.0749063670411985,0,0,0,0
45.617977528089895,67.65188834154351,97.43054488687412,146225,354.9150485436893
0.0749063670411985,0,0,0,0
0.4087193460490463,83.33333333333334,48.62854587773874,424,84.8
0.2724795640326975,75.0,75.52588997769142,79,26.33333333333333
1.7029972752043598,84.0,84.47821405732526,1544,73.52380952380952
2.588555858310627,97.36842105263158,83.21671397873456,3519,95.10810810810811
2.043596730245232,96.66666666666667,72.39061861953026,1958,67.51724137931035
0.4087193460490463,50.0,58.48950453928499,329,109.6666666666667
0.4087193460490463,83.33333333333334,97.17345216621055,193,38.6
0.2724795640326975,100,50.28281516829033,204,51
4.83651226158038,94.36619718309859,83.59701220376321,5773,86.16417910447761
0.4087193460490463,83.33333333333334,65.44227133230056,336,67.2
0.1362397820163488,100,74.80571805931967,257,128.5
0.4087193460490463,83.33333333333334,51.18143943940708,439,87.8
1.7711171662125338,96.15384615384616,67.78225121307666,1364,54.56
4.087193460490464,96.66666666666667,63.62145417767135,4928,84.96551724137932
0.6811989100817438,60.0,53.45854668851414,555,92.5
4.359673024523161,92.1875,71.34519490098421,9650,163.5593220338983
11.10354223433243,92.02453987730061,85.31279837161534,24870,165.8
0.4087193460490463,50.0,90.93472971677974,277,92.33333333333333
0.2724795640326975,100,35.90081422431659,157,39.25
9.673024523160763,99.29577464788733,97.18316933862098,18138,128.6382978723404
1.430517711171662,85.71428571428571,72.3999784864688,1224,68
0.1362397820163488,50.0,36.90833212357519,55,55
52.17983651226158,95.56135770234987,94.57221187379758,127373,174.0068306010929
0.7309941520467835,0,0,0,0
1.900584795321637,7.6923076923076925,1.011042153300025,2,2
0.2923976608187134,50.0,70.01644246718027,1,1
3.654970760233918,16.0,32.67652759350349,12,3
3.8011695906432754,30.76923076923077,25.66482271243108,24,3
1.461988304093567,50.0,9.618284102354172,186,37.2
1.754385964912281,41.66666666666667,17.64391592807164,27,5.4
0.8771929824561403,16.66666666666667,14.60352775632884,19,19
1.608187134502924,18.18181818181818,9.293873452855298,6,3
9.064327485380117,16.12903225806452,19.16584740439515,233,23.3
0.1461988304093567,0,0,0,0
2.339181286549707,25.0,19.59115743012696,15,3.75
1.608187134502924,9.090909090909092,28.03540175623518,5,5
3.362573099415204,26.08695652173913,17.01098205421863,64,10.66666666666667
1.023391812865497,0,0,0,0
6.140350877192982,28.57142857142857,17.91184341737363,236,19.66666666666667
12.57309941520468,25.581395348837212,24.27094322079124,252,11.45454545454546
1.023391812865497,0,0,0,0
0.2923976608187134,0,0,0,0
8.04093567251462,25.454545454545453,19.8720511393036,211,15.07142857142857
2.046783625730994,14.28571428571428,24.80867940583138,89,44.5
1.461988304093567,10.0,53.85819572145232,1,1
33.47953216374269,23.144104803493452,20.26852418461482,825,15.56603773584906
0.1461988304093567,0,0,0,0
1.1695906432748542,0,0,0,0
"@

$lines = $data -split "\r?\n"
$rowIdx = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","

    $ws.Cells.Item($rowIdx, 5).Value  = [double]$parts[0]   # E  particip
    $ws.Cells.Item($rowIdx, 6).Value  = [double]$parts[1]   # F  taxa_sucesso
    $ws.Cells.Item($rowIdx, 12).Value = [double]$parts[2]   # L  apoio_medio
    $ws.Cells.Item($rowIdx, 13).Value = [double]$parts[3]   # M  contribuicoes
    $ws.Cells.Item($rowIdx, 14).Value = [double]$parts[4]   # N  media_contribuicoes

    $rowIdx++
}
